# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) — style matches the other header cells (bold/centered/bordered).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record is constant for every player row (2 through 46): 65 wins, 97 losses, 0 ties.
for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = 65
    $ws.Cells.Item($row, 31).Value = 97
    $ws.Cells.Item($row, 32).Value = 0
}
